# Integrated ELL Comparison figures
# Adds "SpMV Only" (E) and "Comm Only" (F) columns to the first results
# table, removes the stale duplicate second table (rows 13-22), and
# re-sorts the combined table by column A ("nodes").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New column headers (row 2) ---
$ws.Range("E2").Value = "SpMV Only"
$ws.Range("F2").Value = "Comm Only"

# --- New data for column E ("SpMV Only") and F ("Comm Only"), keyed to
#     the existing (p, ppn, nodes) rows in their current (pre-sort)
#     positions 3-11. Row 9 (p=4, ppn=1, nodes=4) has no timing data,
#     same as column D for that row. ---
$ws.Range("E3").Value = 16.529499999999999
$ws.Range("F3").Value = 0.0056

$ws.Range("E4").Value = 10.1092
$ws.Range("F4").Value = 6.8635999999999999

# Row 5 (A5=4, B5=1, C5=4) intentionally left blank in E/F (no data), matching D5.

$ws.Range("E6").Value = 10.4048
$ws.Range("F6").Value = 26.385000000000002

$ws.Range("E7").Value = 5.4991500000000002
$ws.Range("F7").Value = 19.2149

$ws.Range("E8").Value = 2.5464399999999999
$ws.Range("F8").Value = 6.1165300000000002

$ws.Range("E9").Value = 6.5306199999999999
$ws.Range("F9").Value = 28.780999999999999

$ws.Range("E10").Value = 4.0113000000000003
$ws.Range("F10").Value = 7.9760499999999999

$ws.Range("E11").Value = 1.1715
$ws.Range("F11").Value = 7.0250500000000002

# --- Remove the stale duplicate table (rows 13-22) ---
$ws.Rows("13:22").Delete()

# --- Sort the combined table (A3:F11) by column A, ascending ---
$sortObj = $ws.Sort
[void]$sortObj.SortFields.Clear()
[void]$sortObj.SortFields.Add($ws.Range("A3"))
[void]$sortObj.SetRange($ws.Range("A3:F11"))
$sortObj.Header = 0
[void]$sortObj.Apply()

# --- Update selection to match the new extent of the data ---
[void]$ws.Range("F11").Select()
